$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated FilesTab Neo4j query text in cell B4: drop the "File Type" and
# "Breed" columns from the RETURN clause (ICDC Breed 1-14 script fix).
$newQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
         coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newQuery

# The row shrank from 17 to 15 wrapped lines once the two columns were
# removed, so the custom row height comes down from 246.5 to 217.5.
$ws.Rows.Item(4).RowHeight = 217.5

# Reflect the new cursor/selection position left after the edit.
$ws.Range("B4").Select()
